# 19/01/2025 maintenance sekolah & kategori kelas
# Adds a new "Catatan Maintenance DB" note section (rows 29-31) to Sheet1,
# documenting the pembayaran_kelas table auto-creation behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New section: Catatan Maintenance DB -------------------------------
# Header row, highlighted with a yellow fill.
$ws.Range("B29").Value = "Catatan Maintenance DB"
$ws.Range("B29").Interior.Color = 65535

# Detail rows describing the pembayaran_kelas table.
$ws.Range("A30").Value = 1
$ws.Range("B30").Value = "pembayaran_kelas"
$ws.Range("C30").Value = "dibuat otomatis, saat menginput siswa, sistem langsung membuat juga tabel tagihan pembayaran secara otomatis, nantinya admin tinggal melakukan validasi pembayaran "

$ws.Range("C31").Value = "jadi tidak usah ada fitur tambah pembayaran"

# --- View state ----------------------------------------------------------
# Scroll the sheet so row 16 is at the top and select C26, mirroring where
# the author was working when the note was added.
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C26").Select()
